$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.167.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.66%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.800.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.21%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.798.51"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.18%  "

# Row 8
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("E10").Value = "  -1.55%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "

# Row 12
$ws.Range("E12").Value = "  -0.36%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000264"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.53%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.441.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.31%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.798.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.182.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.64%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.44%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.93%  "

# Row 20
$ws.Range("E20").Value = "  -0.23%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.20%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "471.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.61%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.709"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.71%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.19%  "

# Row 25
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000148"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.99%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.75%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.45%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.39%  "

# Row 29
$ws.Range("E29").Value = "  +0.10%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.952.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.21%  "

# Row 31
$ws.Range("E31").Value = "  -2.49%  "

# Row 32
$ws.Range("E32").Value = "  -2.58%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.34%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.757.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.04%  "

# Row 38
$ws.Range("E38").Value = "  -2.13%  "

# Row 39
$ws.Range("E39").Value = "  -7.35%  "

# Row 40
$ws.Range("E40").Value = "  +1.38%  "

# Row 41
$ws.Range("E41").Value = "  +1.20%  "

# Row 42
$ws.Range("E42").Value = "  +0.27%  "

# Row 43
$ws.Range("E43").Value = "  +0.12%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.310"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.29%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.71%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.92%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.53%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "401.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.53%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.33%  "
